$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "**Positive Test Cases:**`n`n**Test Case ID:** 001`n**Description:** User successfully authenticates using OTP sent to official email.`n**Preconditions:** User has access to their official email.`n**Steps:**`n1. User enters their email address on the HR portal.`n2. User requests OTP.`n3. User checks their official email and retrieves the OTP.`n4. User enters the OTP on the HR portal.`n**Expected Results:** User is successfully authenticated and can access HR-related information via WhatsApp.`n**Priority:** High`n`n**Validation Test Cases:**`n`n**Test Case ID:** 002`n**Description:** User enters an incorrect email address.`n**Preconditions:** User has access to the HR portal.`n**Steps:**`n1. User enters an incorrect email address on the HR portal.`n2. User requests OTP.`n**Expected Results:** System displays an error message prompting the user to enter a valid email address.`n**Priority:** Medium`n`n**UI/UX Test Cases:**`n`n**Test Case ID:** 003`n**Description:** OTP request button is clearly visible and accessible on the HR portal.`n**Preconditions:** User has access to the HR portal.`n**Steps:**`n1. User navigates to the authentication page on the HR portal.`n**Expected Results:** OTP request button is prominent and easy to locate on the page.`n**Priority:** Medium`n`n**Performance Test Cases:**`n`n**Test Case ID:** 004`n**Description:** Test system performance when multiple users are simultaneously requesting OTPs.`n**Preconditions:** Multiple users are accessing the HR portal.`n**Steps:**`n1. Multiple users request OTPs simultaneously.`n**Expected Results:** System should be able to handle the load and deliver OTPs promptly to all users.`n**Priority:** High`n`n**Security Test Cases:**`n`n**Test Case ID:** 005`n**Description:** Test if OTPs are valid only for a limited time.`n**Preconditions:** User receives OTP.`n**Steps:**`n1. User retrieves OTP and waits for a significant amount of time.`n2. User tries to enter the same OTP after the validity period has passed.`n**Expected Results:** System should not accept the expired OTP and prompt the user to request a new one.`n**Priority:** High`n`n**Combination Test Cases:**`n`n**Test Case ID:** 006`n**Description:** User enters an incorrect email address and then requests OTP.`n**Preconditions:** User has access to the HR portal.`n**Steps:**`n1. User enters an incorrect email address on the HR portal.`n2. User requests OTP.`n**Expected Results:** System should display an error message for the incorrect email address and not send an OTP.`n**Priority:** Medium"

$ws.Range("B2").Value = $newText
